$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared string text used by D2:D51 (remove the completionTokens field)
$ws.Range("D2:D51").Value = '{"fcCount":1,"fcInfo":{"apiair-conditionerupdate_POST":4}}'

# Move the active selection to D31
$ws.Range("D31").Select()
